$wb = $excel.ActiveWorkbook

# --- Sheet: 2021_13 ---
$ws = $wb.Worksheets.Item("2021_13")
$ws.Range("C4").Value = "'1.3624"
$ws.Range("D4").Value = "'1.332"
$ws.Range("E4").Value = "'1.393"
$ws.Range("C5").Value = "'1.0536"
$ws.Range("D5").Value = "'0.965"
$ws.Range("E5").Value = "'1.150"
$ws.Range("C6").Value = "'1.1105"
$ws.Range("D6").Value = "'1.063"
$ws.Range("E6").Value = "'1.160"
$ws.Range("C7").Value = "'1.6583"
$ws.Range("D7").Value = "'1.606"
$ws.Range("E7").Value = "'1.712"
$ws.Range("C8").Value = "'1.6773"
$ws.Range("D8").Value = "'1.550"
$ws.Range("E8").Value = "'1.814"
$ws.Range("C9").Value = "'1.5101"
$ws.Range("D9").Value = "'1.308"
$ws.Range("E9").Value = "'1.743"
$ws.Range("C10").Value = "'3.2140"
$ws.Range("D10").Value = "'2.522"
$ws.Range("E10").Value = "'4.097"
$ws.Range("C11").Value = "'1.1513"
$ws.Range("D11").Value = "'0.759"
$ws.Range("E11").Value = "'1.746"
$ws.Range("C12").Value = "'2.8159"
$ws.Range("D12").Value = "'1.785"
$ws.Range("E12").Value = "'4.443"
$ws.Range("C13").Value = "'12.8220"
$ws.Range("D13").Value = "'6.029"
$ws.Range("E13").Value = "'27.268"
$ws.Range("C16").Value = "'1.6414"
$ws.Range("D16").Value = "'1.614"
$ws.Range("E16").Value = "'1.669"
$ws.Range("C17").Value = "'1.4450"
$ws.Range("D17").Value = "'1.379"
$ws.Range("E17").Value = "'1.514"
$ws.Range("C18").Value = "'1.6579"
$ws.Range("D18").Value = "'1.614"
$ws.Range("E18").Value = "'1.703"
$ws.Range("C19").Value = "'1.7291"
$ws.Range("D19").Value = "'1.662"
$ws.Range("E19").Value = "'1.798"
$ws.Range("C20").Value = "'1.2950"
$ws.Range("D20").Value = "'1.201"
$ws.Range("E20").Value = "'1.396"
$ws.Range("C21").Value = "'1.3314"
$ws.Range("D21").Value = "'1.150"
$ws.Range("E21").Value = "'1.541"
$ws.Range("C22").Value = "'1.1886"
$ws.Range("D22").Value = "'0.935"
$ws.Range("E22").Value = "'1.511"
$ws.Range("C23").Value = "'0.6578"
$ws.Range("D23").Value = "'0.426"
$ws.Range("E23").Value = "'1.016"
$ws.Range("C24").Value = "'1.3027"
$ws.Range("D24").Value = "'0.776"
$ws.Range("E24").Value = "'2.188"
$ws.Range("C25").Value = "'2.3086"
$ws.Range("D25").Value = "'0.324"
$ws.Range("E25").Value = "'16.470"

# --- Sheet: 2021_20 ---
$ws = $wb.Worksheets.Item("2021_20")
$ws.Range("C4").Value = "'1.2558"
$ws.Range("D4").Value = "'1.232"
$ws.Range("E4").Value = "'1.280"
$ws.Range("C5").Value = "'0.7488"
$ws.Range("D5").Value = "'0.676"
$ws.Range("E5").Value = "'0.829"
$ws.Range("C6").Value = "'1.0962"
$ws.Range("D6").Value = "'1.042"
$ws.Range("E6").Value = "'1.154"
$ws.Range("C7").Value = "'1.1618"
$ws.Range("D7").Value = "'1.122"
$ws.Range("E7").Value = "'1.203"
$ws.Range("C8").Value = "'1.6006"
$ws.Range("D8").Value = "'1.546"
$ws.Range("E8").Value = "'1.657"
$ws.Range("C9").Value = "'1.5308"
$ws.Range("D9").Value = "'1.439"
$ws.Range("E9").Value = "'1.629"
$ws.Range("C10").Value = "'1.3047"
$ws.Range("D10").Value = "'1.159"
$ws.Range("E10").Value = "'1.468"
$ws.Range("C11").Value = "'1.3702"
$ws.Range("D11").Value = "'1.090"
$ws.Range("E11").Value = "'1.722"
$ws.Range("C12").Value = "'1.1949"
$ws.Range("D12").Value = "'0.737"
$ws.Range("E12").Value = "'1.938"
$ws.Range("C13").Value = "'12.2368"
$ws.Range("D13").Value = "'8.343"
$ws.Range("E13").Value = "'17.949"
$ws.Range("C16").Value = "'1.1978"
$ws.Range("D16").Value = "'1.181"
$ws.Range("E16").Value = "'1.215"
$ws.Range("C17").Value = "'1.0309"
$ws.Range("D17").Value = "'0.986"
$ws.Range("E17").Value = "'1.078"
$ws.Range("C18").Value = "'1.1283"
$ws.Range("D18").Value = "'1.099"
$ws.Range("E18").Value = "'1.158"
$ws.Range("C19").Value = "'1.3592"
$ws.Range("D19").Value = "'1.324"
$ws.Range("E19").Value = "'1.396"
$ws.Range("C20").Value = "'1.1468"
$ws.Range("D20").Value = "'1.090"
$ws.Range("E20").Value = "'1.206"
$ws.Range("C21").Value = "'1.3485"
$ws.Range("D21").Value = "'1.228"
$ws.Range("E21").Value = "'1.481"
$ws.Range("C22").Value = "'1.0974"
$ws.Range("D22").Value = "'0.945"
$ws.Range("E22").Value = "'1.274"
$ws.Range("C23").Value = "'2.9842"
$ws.Range("D23").Value = "'2.329"
$ws.Range("E23").Value = "'3.823"
$ws.Range("C24").Value = "'0.9872"
$ws.Range("D24").Value = "'0.707"
$ws.Range("E24").Value = "'1.379"
$ws.Range("C25").Value = "'7.2390"
$ws.Range("D25").Value = "'3.836"
$ws.Range("E25").Value = "'13.660"

# --- Sheet: 2021_24 ---
$ws = $wb.Worksheets.Item("2021_24")
$ws.Range("C4").Value = "'1.0180"
$ws.Range("D4").Value = "'0.996"
$ws.Range("E4").Value = "'1.041"
$ws.Range("C5").Value = "'0.7143"
$ws.Range("D5").Value = "'0.634"
$ws.Range("E5").Value = "'0.804"
$ws.Range("C6").Value = "'0.9900"
$ws.Range("D6").Value = "'0.936"
$ws.Range("E6").Value = "'1.047"
$ws.Range("C7").Value = "'0.9303"
$ws.Range("D7").Value = "'0.890"
$ws.Range("E7").Value = "'0.972"
$ws.Range("C8").Value = "'1.1803"
$ws.Range("D8").Value = "'1.132"
$ws.Range("E8").Value = "'1.230"
$ws.Range("C9").Value = "'1.1977"
$ws.Range("D9").Value = "'1.128"
$ws.Range("E9").Value = "'1.272"
$ws.Range("C10").Value = "'1.3255"
$ws.Range("D10").Value = "'1.211"
$ws.Range("E10").Value = "'1.451"
$ws.Range("C11").Value = "'1.0272"
$ws.Range("D11").Value = "'0.861"
$ws.Range("E11").Value = "'1.226"
$ws.Range("C12").Value = "'3.6197"
$ws.Range("D12").Value = "'2.571"
$ws.Range("E12").Value = "'5.096"
$ws.Range("C13").Value = "'0.6042"
$ws.Range("D13").Value = "'0.310"
$ws.Range("E13").Value = "'1.176"
$ws.Range("C16").Value = "'1.1145"
$ws.Range("D16").Value = "'1.099"
$ws.Range("E16").Value = "'1.130"
$ws.Range("C17").Value = "'0.9437"
$ws.Range("D17").Value = "'0.902"
$ws.Range("E17").Value = "'0.987"
$ws.Range("C18").Value = "'1.0595"
$ws.Range("D18").Value = "'1.032"
$ws.Range("E18").Value = "'1.088"
$ws.Range("C19").Value = "'1.1438"
$ws.Range("D19").Value = "'1.116"
$ws.Range("E19").Value = "'1.173"
$ws.Range("C20").Value = "'1.3989"
$ws.Range("D20").Value = "'1.351"
$ws.Range("E20").Value = "'1.448"
$ws.Range("C21").Value = "'1.1421"
$ws.Range("D21").Value = "'1.065"
$ws.Range("E21").Value = "'1.225"
$ws.Range("C22").Value = "'1.1921"
$ws.Range("D22").Value = "'1.064"
$ws.Range("E22").Value = "'1.336"
$ws.Range("C23").Value = "'1.3888"
$ws.Range("D23").Value = "'1.133"
$ws.Range("E23").Value = "'1.703"
$ws.Range("C24").Value = "'2.4392"
$ws.Range("D24").Value = "'1.808"
$ws.Range("E24").Value = "'3.291"
$ws.Range("C25").Value = "'1.4128"
$ws.Range("D25").Value = "'0.986"
$ws.Range("E25").Value = "'2.025"

# --- Sheet: 2021_30 ---
$ws = $wb.Worksheets.Item("2021_30")
$ws.Range("C4").Value = "'0.6885"
$ws.Range("D4").Value = "'0.664"
$ws.Range("E4").Value = "'0.714"
$ws.Range("C5").Value = "'0.4804"
$ws.Range("D5").Value = "'0.400"
$ws.Range("E5").Value = "'0.577"
$ws.Range("C6").Value = "'0.6696"
$ws.Range("D6").Value = "'0.617"
$ws.Range("E6").Value = "'0.727"
$ws.Range("C7").Value = "'0.6393"
$ws.Range("D7").Value = "'0.596"
$ws.Range("E7").Value = "'0.686"
$ws.Range("C8").Value = "'0.6676"
$ws.Range("D8").Value = "'0.619"
$ws.Range("E8").Value = "'0.720"
$ws.Range("C9").Value = "'1.0305"
$ws.Range("D9").Value = "'0.922"
$ws.Range("E9").Value = "'1.151"
$ws.Range("C10").Value = "'1.1486"
$ws.Range("D10").Value = "'0.991"
$ws.Range("E10").Value = "'1.332"
$ws.Range("C11").Value = "'1.3185"
$ws.Range("D11").Value = "'1.080"
$ws.Range("E11").Value = "'1.610"
$ws.Range("C12").Value = "'1.1258"
$ws.Range("D12").Value = "'0.818"
$ws.Range("E12").Value = "'1.549"
$ws.Range("C13").Value = "'1.2602"
$ws.Range("D13").Value = "'0.853"
$ws.Range("E13").Value = "'1.861"
$ws.Range("C16").Value = "'0.9658"
$ws.Range("D16").Value = "'0.953"
$ws.Range("E16").Value = "'0.978"
$ws.Range("C17").Value = "'0.8414"
$ws.Range("D17").Value = "'0.805"
$ws.Range("E17").Value = "'0.880"
$ws.Range("C18").Value = "'0.8919"
$ws.Range("D18").Value = "'0.869"
$ws.Range("E18").Value = "'0.915"
$ws.Range("C19").Value = "'0.9701"
$ws.Range("D19").Value = "'0.947"
$ws.Range("E19").Value = "'0.994"
$ws.Range("C20").Value = "'0.9935"
$ws.Range("D20").Value = "'0.964"
$ws.Range("E20").Value = "'1.024"
$ws.Range("C21").Value = "'1.3713"
$ws.Range("D21").Value = "'1.303"
$ws.Range("E21").Value = "'1.443"
$ws.Range("C22").Value = "'1.1921"
$ws.Range("D22").Value = "'1.102"
$ws.Range("E22").Value = "'1.289"
$ws.Range("C23").Value = "'1.7743"
$ws.Range("D23").Value = "'1.535"
$ws.Range("E23").Value = "'2.051"
$ws.Range("C24").Value = "'1.2497"
$ws.Range("D24").Value = "'0.975"
$ws.Range("E24").Value = "'1.602"
$ws.Range("C25").Value = "'0.3086"
$ws.Range("D25").Value = "'0.197"
$ws.Range("E25").Value = "'0.484"
